$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(4)
$p.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item(4)
$r = $newPara.Range
$r.Text = "OBJECTIVE: Eager to drive back-end solutions at US Synthetic on a full-time basis"

$p2 = $d.Paragraphs.Item(4)
$prefixLen = "OBJECTIVE: ".Length
$afterRange = $d.Range($p2.Range.Start + $prefixLen, $p2.Range.End - 1)
$afterRange.Font.Size = 12

# select from mark (End-1) through start of next paragraph's first char
$p5 = $d.Paragraphs.Item(5)
$spanRange = $d.Range($p2.Range.End - 1, $p5.Range.Start + 1)
Write-Output "span: [$($spanRange.Text)]"
$spanRange.Font.Size = 12

$p3 = $d.Paragraphs.Item(4)
Write-Output "final p4: [$($p3.Range.Text)]"
$p6 = $d.Paragraphs.Item(5)
Write-Output "final p5: [$($p6.Range.Text)]"
